$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fill in the Observed Result (I) and Pass/Fail (J) columns for rows
# that were previously left blank, matching the Expected Result (H)
# and marking them as Pass - reflecting completed test cases.

$ws.Range("I4").Value = $ws.Range("H4").Value2
$ws.Range("J4").Value = "Pass"

$ws.Range("I5").Value = $ws.Range("H5").Value2
$ws.Range("J5").Value = "Pass"

$ws.Range("I9").Value = $ws.Range("H9").Value2
$ws.Range("J9").Value = "Pass"

# Row 9 ("Delete an expense from list") now wraps across three lines once the
# Observed Result column is populated, so the row grows taller to fit.
$ws.Rows.Item(9).RowHeight = 57.6

# Update the view to reflect scrolling to the newly-completed row (TS_03 / "Delete an expense from list")
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G9").Select()
